# Insert a new data row at row 129, shifting existing rows 129:201 down to 130:202.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("129").Insert()

# Populate the newly inserted row 129 with the new record's data.
$ws.Range("A129").Value = 1
$ws.Range("B129").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C129").Value = "Arica y Parinacota"
$ws.Range("D129").Value = 44488
$ws.Range("E129").Value = 15
$ws.Range("F129").Value = 100114013
$ws.Range("G129").Value = "Zanahoria"
$ws.Range("H129").Value = "Sin especificar"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 70
$ws.Range("K129").Value = 11000
$ws.Range("L129").Value = 12000
$ws.Range("M129").Value = 11500
$ws.Range("N129").Value = "`$/saco 25 kilos"
$ws.Range("O129").Value = "Región de Arica y Parinacota"
$ws.Range("P129").Value = 460
$ws.Range("Q129").Value = 25
$ws.Range("R129").Value = "Hortaliza"
